$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 8: style change only (fill becomes the new light theme fill) ---
$ws.Range("A8:E8").Interior.ThemeColor = 3

# --- Row 11: content + style changes ---
$ws.Range("C11").Value2 = "Mario"
$ws.Range("B26").Value2 = "Reviews"
$ws.Range("E11").Value2 = "Implemented as DAO methods"
$ws.Range("B21").Value2 = "Rate products"

$ws.Range("D11").Value2 = "Finished"
$ws.Range("A11:E11").Interior.ThemeColor = 3

# --- Row 12: content + style changes ---
$ws.Range("D12").Value2 = "Finished"
$ws.Range("A12:E12").Interior.ThemeColor = 3

# --- Row 13: content + style changes ---
$ws.Range("D13").Value2 = "Finished"
$ws.Range("A13:E13").Interior.ThemeColor = 3

# --- Row 21: new task row (Mandatory Tasks table), style unchanged ---
$ws.Range("A21").Value2 = 17
$ws.Range("C21").Value2 = "Not Assigned"
$ws.Range("D21").Value2 = "Not Started"

# --- Row 26: new task row (Bonus Tasks table), style unchanged ---
$ws.Range("C26").Value2 = "Not Assigned"
$ws.Range("D26").Value2 = "Not Started"

# --- Update view: scrolled down & new selection ---
$ws.Activate()
$ws.Range("E21").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
